$d = $word.ActiveDocument

# The document currently ends with a paragraph that holds an inline
# drawing (the last paragraph of the body, right before the final
# sectPr). The diff appends four new plain-text paragraphs right after
# that drawing paragraph.
#
# A plain InsertParagraphAfter() on that paragraph's range would make
# the new paragraph mark inherit the drawing run's rPr (w:noProof),
# which the target XML does not have. To avoid that, we capture the
# drawing paragraph's own OOXML verbatim and re-insert it together with
# the new plain paragraphs via Range.InsertXML, which lets us control
# the exact markup of what gets written (no inherited formatting on the
# new runs).

$picPara = $d.Paragraphs.Last

# Pull this paragraph's real WordProcessingML out of the flat-OPC
# payload that WordOpenXML returns, keeping it byte-for-byte so the
# drawing/relationship (e.g. r:embed) stays intact.
$wordXml = $picPara.Range.WordOpenXML
if (-not ($wordXml -match '(?s)<w:body>(<w:p .*?</w:p>)\s*(?:<w:p[ />]|<w:sectPr)')) {
    throw "Could not locate the drawing paragraph's OOXML"
}
$picXml = $matches[1]

function Escape-XmlText([string]$text) {
    $text = $text -replace '&', '&amp;'
    $text = $text -replace '<', '&lt;'
    $text = $text -replace '>', '&gt;'
    return $text
}

$newParagraphTexts = @(
    "Para que el disparo sea efectivo contra l disparo de el cañon 1",
    "X1 = X2",
    "Xo1+Vx1*T1 = Xo2 + Vx2*T2",
    "(XD,YD) <= 0,025d(XO,YO)"
)

$wNs = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"
$newXml = ""
foreach ($t in $newParagraphTexts) {
    $escaped = Escape-XmlText $t
    $newXml += "<w:p xmlns:w=`"$wNs`"><w:r><w:t>$escaped</w:t></w:r></w:p>"
}

$r = $picPara.Range
$r.Collapse(0)   # wdCollapseEnd
$r.InsertXML($picXml + $newXml) | Out-Null
